$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.477.08'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '3.412.07'
$ws.Range('E3').Value = '  +4.08%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '653.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.48'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.423'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.00%  '
$ws.Range('E9').Value = '  +9.56%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.407.56'
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.212'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +14.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000259'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.25%  '
$ws.Range('D16').Value = '97.200.36'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '4.038.98'
$ws.Range('E17').Value = '  +4.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +37.11%  '
$ws.Range('D19').Value = '3.396.11'
$ws.Range('E19').Value = '  +3.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +14.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.516'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +61.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +17.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '506.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000205'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '98.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.55%  '
$ws.Range('D29').Value = '3.606.62'
$ws.Range('E29').Value = '  +5.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.155'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +14.69%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.199'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.46%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.574'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +21.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '29.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +17.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.154'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '513.95'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.857'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.41%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0422'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +25.79%  '
$ws.Range('B45').Value = 'MantraDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +15.80%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +16.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.02%  '
